$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.922.28"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.813.49"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.28"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4643"
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3656"
$ws.Range("E8").Value = "  -1.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07359"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8680"
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.832.21"
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.361"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07104"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.499"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.05"
$ws.Range("E16").Value = "  -1.51%  "
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008698"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.61"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.933.76"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.283"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.58"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.047.60"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.80"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.30"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.117"
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.249"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.52"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08905"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7540"
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.161"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.903"
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.086"
$ws.Range("E37").Value = "  -1.66%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05289"
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("E39").Value = "  -1.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.980"
$ws.Range("E40").Value = "  +1.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.245"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5290"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.282"
$ws.Range("E43").Value = "  -5.63%  "
$ws.Range("E44").Value = "  -0.72%  "
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4854"
$ws.Range("E46").Value = "  -2.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.39"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("E48").Value = "  +0.20%  "
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.16"
$ws.Range("E50").Value = "  -0.76%  "
$ws.Range("E51").Value = "  -0.04%  "
